# "add more filter examples"
# Renames the existing sheet to "simple" and adds a new "monthly report"
# worksheet containing a small sales/profit table with an AutoFilter,
# a color-scale conditional format, a data validation list and a quarterly
# running-total column, mirroring the target workbook.

$wb = $excel.ActiveWorkbook

# --- rename original sheet -------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "simple"

# --- add the new worksheet right after "simple" -----------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "monthly report"

# --- title row ---------------------------------------------------------------
$ws2.Range("A1").Value = "Cape Company Monthly Report"
$ws2.Range("A1:F1").Merge()
$ws2.Rows.Item(1).RowHeight = 24

# --- header row ---------------------------------------------------------------
$ws2.Range("A3").Value = "Date"
$ws2.Range("B3").Value = "Customer No."
$ws2.Range("C3").Value = "Income"
$ws2.Range("D3").Value = "Cost"
$ws2.Range("E3").Value = "Monthly Profit"
$ws2.Range("F3").Value = "Quarterly Profit"
$ws2.Rows.Item(3).RowHeight = 21

# --- data rows -----------------------------------------------------------------
$dates = 42384,42415,42444,42475,42505,42536,42566,42597,42628,42658,42689,42719
$custNo = 13,19,25,22,28,35,20,31,27,24,19,17
$income = 89500,100500,119200,115900,123700,129300,110700,125100,120100,118400,100300,94200
$cost   = 62600,60300,27800,79600,84000,83100,77300,85500,78900,91000,65100,65800

for ($i = 0; $i -lt 12; $i++) {
    $r = 4 + $i
    $ws2.Cells.Item($r, 1).Value = $dates[$i]
    $ws2.Cells.Item($r, 2).Value = $custNo[$i]
    $ws2.Cells.Item($r, 3).Value = $income[$i]
    $ws2.Cells.Item($r, 4).Value = $cost[$i]
    $ws2.Cells.Item($r, 5).Formula = "=C" + $r + "-D" + $r
    $ws2.Rows.Item($r).RowHeight = 18
}

# quarterly running total in column F (resets every 3 rows)
$ws2.Range("F4").Formula = "=E4"
$ws2.Range("F5").Formula = "=E4+E5"
$ws2.Range("F6").Formula = "=E4+E5+E6"
$ws2.Range("F7").Formula = "=E7"
$ws2.Range("F8").Formula = "=E7+E8"
$ws2.Range("F9").Formula = "=E7+E8+E9"
$ws2.Range("F10").Formula = "=E10"
$ws2.Range("F11").Formula = "=E10+E11"
$ws2.Range("F12").Formula = "=E10+E11+E12"
$ws2.Range("F13").Formula = "=E13"
$ws2.Range("F14").Formula = "=E13+E14"
$ws2.Range("F15").Formula = "=E13+E14+E15"

# --- number formats -------------------------------------------------------------
$ws2.Range("A4:A15").NumberFormat = "[$-409]d-mmm;@"
$ws2.Range("C4:F15").NumberFormat = "$#,##0"

# --- column widths ---------------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 11.1640625
$ws2.Columns.Item(2).ColumnWidth = 15
$ws2.Columns.Item(3).ColumnWidth = 12.1640625
$ws2.Columns.Item(4).ColumnWidth = 10
$ws2.Columns.Item(5).ColumnWidth = 18.83203125
$ws2.Columns.Item(6).ColumnWidth = 18.1640625

# --- table look & feel ------------------------------------------------------------
$ws2.Range("A3:F15").AutoFilter(1)

$cf = $ws2.Range("C4:C15").FormatConditions.AddColorScale(2)
$cf.ColorScaleCriteria.Item(1).Type = 1
$cf.ColorScaleCriteria.Item(1).FormatColor.Color = 10285055
$cf.ColorScaleCriteria.Item(2).Type = 2
$cf.ColorScaleCriteria.Item(2).FormatColor.Color = 2650623

$ws2.Range("L4:L11").Validation.Add(3, 1, 1, "=`$N`$11")

$ws2.PageSetup.Orientation = 1

# --- selection bookkeeping (keep "simple" as the active tab) ----------------------
$ws2.Range("F4").Select()
$ws1.Activate()
